$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: create additional blank styled buffer rows (1254:1275) by copying format from row 1253 ---
$ws.Range("A1253:J1253").Copy() | Out-Null
$ws.Range("A1254:J1275").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Step 2: update D1236 (effort value correction) ---
$ws.Range("D1236").Value = '0h 45m'

# --- Step 3: fill in new log rows 1243:1258 ---
$ws.Range("A1243").Value = '2024-10-27'
$ws.Range("B1243").Value = '10:00'
$ws.Range("C1243").Value = '12:45'
$ws.Range("D1243").Value = '2h 45m'
$ws.Range("E1243").Value = '#studying'
$ws.Range("F1243").Value = 'NB20230816, 58 -> 94 (36 pp).'
$ws.Range("G1243").Value = '''False'
$ws.Range("H1243").Value = '''False'

$ws.Range("A1244").Value = '2024-10-27'
$ws.Range("B1244").Value = '14:15'
$ws.Range("C1244").Value = '17:00'
$ws.Range("D1244").Value = '2h 45m'
$ws.Range("E1244").Value = '#studying'
$ws.Range("F1244").Value = 'NB20230816, 58 -> 94 (36 pp).'
$ws.Range("G1244").Value = '''False'
$ws.Range("H1244").Value = '''False'

$ws.Range("A1245").Value = '2024-10-27'
$ws.Range("B1245").Value = '19:15'
$ws.Range("C1245").Value = '20:00'
$ws.Range("D1245").Value = '0h 45m'
$ws.Range("E1245").Value = '#studying'
$ws.Range("F1245").Value = 'NB20230816, 58 -> 94 (36 pp).'
$ws.Range("G1245").Value = '''False'
$ws.Range("H1245").Value = '''False'

$ws.Range("A1246").Value = '2024-10-27'
$ws.Range("B1246").Value = '20:30'
$ws.Range("C1246").Value = '22:00'
$ws.Range("D1246").Value = '1h 30m'
$ws.Range("E1246").Value = '#studying'
$ws.Range("F1246").Value = 'NB20230816, 58 -> 94 (36 pp).'
$ws.Range("G1246").Value = '''False'
$ws.Range("H1246").Value = '''False'

$ws.Range("A1247").Value = '2024-10-28'
$ws.Range("B1247").Value = '10:15'
$ws.Range("C1247").Value = '14:00'
$ws.Range("D1247").Value = '3h 45m'
$ws.Range("E1247").Value = '#studying'
$ws.Range("F1247").Value = 'NB20230816, 95 -> 116 (21 pp).'
$ws.Range("G1247").Value = '''False'
$ws.Range("H1247").Value = '''False'

$ws.Range("A1248").Value = '2024-10-28'
$ws.Range("B1248").Value = '15:30'
$ws.Range("C1248").Value = '16:15'
$ws.Range("D1248").Value = '0h 45m'
$ws.Range("E1248").Value = '#python'
$ws.Range("F1248").Value = 'nwshared v1.6.0'
$ws.Range("G1248").Value = '''True'
$ws.Range("H1248").Value = '''True'

$ws.Range("A1249").Value = '2024-10-28'
$ws.Range("B1249").Value = '16:15'
$ws.Range("C1249").Value = '17:45'
$ws.Range("D1249").Value = '1h 30m'
$ws.Range("E1249").Value = '#python'
$ws.Range("F1249").Value = 'nwreadinglist v3.8.0'
$ws.Range("G1249").Value = '''True'
$ws.Range("H1249").Value = '''True'

$ws.Range("A1250").Value = '2024-10-28'
$ws.Range("B1250").Value = '17:45'
$ws.Range("C1250").Value = '18:45'
$ws.Range("D1250").Value = '1h 00m'
$ws.Range("E1250").Value = '#python'
$ws.Range("F1250").Value = 'nwtimetracking v3.8.0'
$ws.Range("G1250").Value = '''True'
$ws.Range("H1250").Value = '''True'

$ws.Range("A1251").Value = '2024-10-28'
$ws.Range("B1251").Value = '20:45'
$ws.Range("C1251").Value = '21:15'
$ws.Range("D1251").Value = '0h 30m'
$ws.Range("E1251").Value = '#python'
$ws.Range("F1251").Value = 'nwtraderaanalytics v4.3.0'
$ws.Range("G1251").Value = '''True'
$ws.Range("H1251").Value = '''True'

$ws.Range("A1252").Value = '2024-10-29'
$ws.Range("B1252").Value = '10:00'
$ws.Range("C1252").Value = '15:30'
$ws.Range("D1252").Value = '5h 30m'
$ws.Range("E1252").Value = '#python'
$ws.Range("F1252").Value = 'nwpackageversions v1.6.0'
$ws.Range("G1252").Value = '''True'
$ws.Range("H1252").Value = '''False'

$ws.Range("A1253").Value = '2024-10-31'
$ws.Range("B1253").Value = '08:00'
$ws.Range("C1253").Value = '08:45'
$ws.Range("D1253").Value = '0h 45m'
$ws.Range("E1253").Value = '#python'
$ws.Range("F1253").Value = 'nwpackageversions v1.6.0'
$ws.Range("G1253").Value = '''True'
$ws.Range("H1253").Value = '''False'

$ws.Range("A1254").Value = '2024-10-31'
$ws.Range("B1254").Value = '17:00'
$ws.Range("C1254").Value = '17:30'
$ws.Range("D1254").Value = '0h 30m'
$ws.Range("E1254").Value = '#python'
$ws.Range("F1254").Value = 'nwpackageversions v1.6.0'
$ws.Range("G1254").Value = '''True'
$ws.Range("H1254").Value = '''False'

$ws.Range("A1255").Value = '2024-11-01'
$ws.Range("B1255").Value = '08:00'
$ws.Range("C1255").Value = '08:45'
$ws.Range("D1255").Value = '0h 45m'
$ws.Range("E1255").Value = '#python'
$ws.Range("F1255").Value = 'nwpackageversions v1.6.0'
$ws.Range("G1255").Value = '''True'
$ws.Range("H1255").Value = '''False'

$ws.Range("A1256").Value = '2024-11-01'
$ws.Range("B1256").Value = '17:00'
$ws.Range("C1256").Value = '17:30'
$ws.Range("D1256").Value = '0h 30m'
$ws.Range("E1256").Value = '#python'
$ws.Range("F1256").Value = 'nwpackageversions v1.6.0'
$ws.Range("G1256").Value = '''True'
$ws.Range("H1256").Value = '''False'

$ws.Range("A1257").Value = '2024-11-03'
$ws.Range("B1257").Value = '12:15'
$ws.Range("C1257").Value = '14:30'
$ws.Range("D1257").Value = '2h 15m'
$ws.Range("E1257").Value = '#python'
$ws.Range("F1257").Value = 'nwpackageversions v1.6.0'
$ws.Range("G1257").Value = '''True'
$ws.Range("H1257").Value = '''False'

$ws.Range("A1258").Value = '2024-11-03'
$ws.Range("B1258").Value = '15:30'
$ws.Range("C1258").Value = '17:30'
$ws.Range("D1258").Value = '2h 00m'
$ws.Range("E1258").Value = '#python'
$ws.Range("F1258").Value = 'nwpackageversions v1.6.0'
$ws.Range("G1258").Value = '''True'
$ws.Range("H1258").Value = '''False'

# --- Step 4: fill in I/J formulas, matching the original shared-formula batch groupings ---
$ws.Range("I1243").Formula = "=YEAR(A1243)"
$ws.Range("J1243").Formula = "=MONTH(A1243)"
$ws.Range("I1244:I1246").Formula = "=YEAR(A1244)"
$ws.Range("J1244:J1246").Formula = "=MONTH(A1244)"
$ws.Range("I1247:I1248").Formula = "=YEAR(A1247)"
$ws.Range("J1247:J1248").Formula = "=MONTH(A1247)"
$ws.Range("I1249").Formula = "=YEAR(A1249)"
$ws.Range("J1249").Formula = "=MONTH(A1249)"
$ws.Range("I1250").Formula = "=YEAR(A1250)"
$ws.Range("J1250").Formula = "=MONTH(A1250)"
$ws.Range("I1251").Formula = "=YEAR(A1251)"
$ws.Range("J1251").Formula = "=MONTH(A1251)"
$ws.Range("I1252").Formula = "=YEAR(A1252)"
$ws.Range("J1252").Formula = "=MONTH(A1252)"
$ws.Range("I1253:I1254").Formula = "=YEAR(A1253)"
$ws.Range("J1253:J1254").Formula = "=MONTH(A1253)"
$ws.Range("I1255:I1256").Formula = "=YEAR(A1255)"
$ws.Range("J1255:J1256").Formula = "=MONTH(A1255)"
$ws.Range("I1257:I1258").Formula = "=YEAR(A1257)"
$ws.Range("J1257:J1258").Formula = "=MONTH(A1257)"

# --- Step 5: update sheet view (frozen pane + selection) to match final state ---
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D1260").Select() | Out-Null
